# ============================================================================
# Edit script: updates the Dheeraj Chand resume (topographic_classic, short)
# to match the target revision described in the task's unified diff.
# ============================================================================

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper functions
# ---------------------------------------------------------------------------

function Get-ParaIndex {
    param($text)
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.TrimEnd() -eq $text) {
            return $i
        }
    }
    return -1
}

function Set-ParaText {
    param($index, $text)
    $d.Paragraphs.Item($index).Range.Text = $text
}

function Insert-ParaAfter {
    # Inserts a new paragraph right after $afterIndex, sets its text/style,
    # and returns the new paragraph's index (always $afterIndex + 1).
    param($afterIndex, $text, $style)
    $anchor = $d.Paragraphs.Item($afterIndex).Range
    $anchor.InsertParagraphAfter()
    $newIndex = $afterIndex + 1
    $np = $d.Paragraphs.Item($newIndex)
    $np.Range.Text = $text
    if ($style) {
        $np.Style = $style
    }
    return $newIndex
}

function Delete-ParaRange {
    # Deletes paragraphs fromIndex..toIndex inclusive (1-based, by index
    # computed BEFORE the delete happens).
    param($fromIndex, $toIndex)
    $start = $d.Paragraphs.Item($fromIndex).Range.Start
    $end = $d.Paragraphs.Item($toIndex).Range.End
    $range = $d.Range($start, $end)
    $range.Delete()
}

Write-Host "Initial paragraph count: $($d.Paragraphs.Count)"

# ---------------------------------------------------------------------------
# 1. Remove the contact-info paragraph entirely (phone/email/site/linkedin).
# ---------------------------------------------------------------------------
$contactIdx = Get-ParaIndex "+1 (512) 555-0123 | dheeraj@dheerajchand.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX"
if ($contactIdx -ne -1) {
    $d.Paragraphs.Item($contactIdx).Range.Delete()
}

Write-Host "After contact removal: $($d.Paragraphs.Count)"

# ---------------------------------------------------------------------------
# 2. Professional summary: swap in the new, metrics-forward summary text.
# ---------------------------------------------------------------------------
$summaryIdx = Get-ParaIndex "Experienced data scientist and software engineer with 15+ years of expertise in geospatial analysis, demographic research, and political data. Proven track record of building scalable systems, conducting complex analyses, and delivering actionable insights for campaigns, organizations, and government agencies."
Set-ParaText $summaryIdx "Senior data scientist and software engineer specializing in geospatial machine learning and large-scale demographic analysis. Developed algorithms that improved demographic classification accuracy from 23% to 64%, processed data across 178,000+ precincts, and built platforms serving thousands of analysts nationwide."

# ---------------------------------------------------------------------------
# 3. Core competencies: blank out the competency-tag line (becomes an empty
#    paragraph, heading stays).
# ---------------------------------------------------------------------------
$coreIdx = Get-ParaIndex "CODE • COMPUTE • INTERACT • MEASURE • PLATFORMS • TRACK"
Set-ParaText $coreIdx ""

Write-Host "After summary/core edits: $($d.Paragraphs.Count)"

# ---------------------------------------------------------------------------
# 4. Siege Analytics (Partner) role: new dates, new subtitle, new bullets.
# ---------------------------------------------------------------------------
$siegeIdx = Get-ParaIndex "Partner - Siege Analytics (Austin, TX) | 2020 - Present"
Set-ParaText $siegeIdx "Partner - Siege Analytics (Austin, TX) | 2005 - Present"

$siegeSubIdx = Get-ParaIndex "Data Science & Political Analytics"
Set-ParaText $siegeSubIdx "Data, Technology and Strategy Consulting"

$siegeBullet1Idx = Get-ParaIndex "• Uncovered decades of demographic miscoding in voter files, discovering 2.7M previously mischaracterized Democratic voters"
Set-ParaText $siegeBullet1Idx "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%"

$siegeBullet2Idx = Get-ParaIndex "• Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States"
Set-ParaText $siegeBullet2Idx "• Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration"

$siegeBullet3Idx = Get-ParaIndex "• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct redistricting analysis"
Set-ParaText $siegeBullet3Idx "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%"

Write-Host "After Siege Analytics edits: $($d.Paragraphs.Count)"

# ---------------------------------------------------------------------------
# 5. Professional experience: replace the five legacy roles (Lake Research
#    Partners/Senior Data Scientist through The Feldman Group) with the
#    seven roles in the new, fuller work history.
# ---------------------------------------------------------------------------
$oldJobsStart = Get-ParaIndex "Senior Data Scientist - Lake Research Partners (Washington, DC) | 2018 - 2020"
$oldJobsEnd = Get-ParaIndex "• Trained staff on PHP/MySQL for data analysis and reporting systems"
Delete-ParaRange $oldJobsStart $oldJobsEnd

Write-Host "After legacy job-history removal: $($d.Paragraphs.Count)"

# Anchor: insert the new job history right after the last Siege Analytics bullet.
$anchorIdx = Get-ParaIndex "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%"

$idx = $anchorIdx
$idx = Insert-ParaAfter $idx "Data Products Manager - Helm/Murmuration (Austin, TX) | June 2021 - May 2023" "Heading3"
$idx = Insert-ParaAfter $idx "Civic Graph & Civic Pulse Director" $null
$idx = Insert-ParaAfter $idx "• Conceived, architected and built Civic Graph multi-tenant data warehouse processing government data from Census, Bureau of Labor Statistics, National Council of Educational Statistics" $null
$idx = Insert-ParaAfter $idx "• Built multi-dimensional data warehouse measuring socio-economic changes in America at every level across attitudinal, behavioral, demographic, economic and geographical dimensions" $null
$idx = Insert-ParaAfter $idx "• Managed engineering teams of 7-11 professionals while setting technical direction for data architecture" $null

$idx = Insert-ParaAfter $idx "Analytics Supervisor - GSD&M (Austin, TX) | November 2019 - June 2020" "Heading3"
$idx = Insert-ParaAfter $idx "Big Data Engineering Transformation" $null
$idx = Insert-ParaAfter $idx "• Transformed small data team into big data engineering team, scaling from laptop datasets to Hadoop Clusters and Hive on AWS" $null
$idx = Insert-ParaAfter $idx "• Managed accounts including United States Air Force, Southwest Airlines/Chase and Indeed" $null
$idx = Insert-ParaAfter $idx "• Rewrote mission and offerings of department and drafted integration plan with strategy team" $null

$idx = Insert-ParaAfter $idx "Software Engineer - Mautinoa Technologies (Austin, TX) | August 2016 - February 2018" "Heading3"
$idx = Insert-ParaAfter $idx "SimCrisis Product Owner/Engineer" $null
$idx = Insert-ParaAfter $idx "• Conceived, architected and engineered econometric simulation software for humanitarian crises intervention measurement" $null
$idx = Insert-ParaAfter $idx "• Built SimCrisis GeoDjango web application using multi-agent modeling to create econometric simulations of crisis economies" $null
$idx = Insert-ParaAfter $idx "• Designed modular application accepting rules extensions for ethnic strife, different crises/disasters, supply failures" $null

$idx = Insert-ParaAfter $idx "Senior Analyst - Myers Research (Austin, TX) | August 2012 - February 2014" "Heading3"
$idx = Insert-ParaAfter $idx "RACSO Product Owner/Engineer" $null
$idx = Insert-ParaAfter $idx "• Designed comprehensive survey instruments for specialized voting segments and niche markets" $null
$idx = Insert-ParaAfter $idx "• Co-developed RACSO web application managing all aspects of survey operations from instrument design to data analysis" $null
$idx = Insert-ParaAfter $idx "• Wrote RFP and analyzed bids from 1,200 vendors for research platform development" $null

$idx = Insert-ParaAfter $idx "Research Director - PCCC (Washington, DC) | 2010 - 2012" "Heading3"
$idx = Insert-ParaAfter $idx "Political Research & Data Analysis (FLEEM System)" $null
$idx = Insert-ParaAfter $idx "• Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of simultaneous phone calls using emulated predictive dialer for regulated political surveys" $null
$idx = Insert-ParaAfter $idx "• Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren" $null
$idx = Insert-ParaAfter $idx "• Built comprehensive tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver" $null

$idx = Insert-ParaAfter $idx "Software Engineer - Salsa Labs (Washington, DC) | January 2011 - August 2011" "Heading3"
$idx = Insert-ParaAfter $idx "Geospatial CRM Development" $null
$idx = Insert-ParaAfter $idx "• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands simultaneously" $null
$idx = Insert-ParaAfter $idx "• Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers" $null
$idx = Insert-ParaAfter $idx "• Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill" $null

$idx = Insert-ParaAfter $idx "Programmer - Lake Research Partners (Washington, DC) | April 2008 - December 2008" "Heading3"
$idx = Insert-ParaAfter $idx "Political Research & Analytics" $null
$idx = Insert-ParaAfter $idx "• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party" $null
$idx = Insert-ParaAfter $idx "• Harmonized data from 20+ polling firms with incompatible methodologies and encoding systems" $null
$idx = Insert-ParaAfter $idx "• Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+" $null

Write-Host "After new job-history insertion: $($d.Paragraphs.Count)"
